$d = $word.ActiveDocument

function Insert-XmlFragment($range, $innerBodyXml) {
    $xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$innerBodyXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $range.InsertXML($xml)
}

# --- Step 1: replace the 2nd (empty) paragraph with the new title / abstract block ---
$target = $d.Paragraphs.Item(2)

$newBlock = @'
<w:p>
  <w:pPr>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:u w:val="single"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t>Hyperloop On-Boarding User Guide</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:u w:val="single"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pBdr>
      <w:bottom w:val="single" w:sz="12" w:space="1" w:color="auto"/>
    </w:pBdr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">Abstract: </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">The purpose of this document is to provide new users with a guide to bring up the Hyperloop system </w:t>
  </w:r>
  <w:r>
    <w:t>so they can be active contributors to the project.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pBdr>
      <w:bottom w:val="single" w:sz="12" w:space="1" w:color="auto"/>
    </w:pBdr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="center"/>
  </w:pPr>
</w:p>
'@

Insert-XmlFragment $target.Range $newBlock

# --- Step 2: rewrite the "Software : version : software user-guide" paragraph ---
# (it was pushed down by 4 extra paragraphs, so it is now paragraph #7)
$softwarePara = $d.Paragraphs.Item(7)

$softwareBlock = @'
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Software : version</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve"> : software user-guide</w:t>
  </w:r>
</w:p>
'@

Insert-XmlFragment $softwarePara.Range $softwareBlock

# --- Step 3: rewrite the "Hardware : Part Number : manufacturer user-guide" paragraph ---
$hwIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Hardware*Part Number*") {
        $hwIndex = $i
        break
    }
}
$hardwarePara = $d.Paragraphs.Item($hwIndex)

$hardwareBlock = @'
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Hardware</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve"> : Part Number</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve"> : manufacturer user-guide</w:t>
  </w:r>
</w:p>
'@

Insert-XmlFragment $hardwarePara.Range $hardwareBlock

# --- Step 4: add a new empty paragraph after the very last paragraph of the body ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

Write-Output "edit complete"
